$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.600.09"

$ws.Range("D3").Value = "1.693.49"
$ws.Range("E3").Value = "  -5.67%  "

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.006"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.27%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "219.70"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -5.00%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.5098"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -13.19%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "1.007"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.21%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.2651"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -3.97%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "22.13"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -4.35%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.06333"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -6.02%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.07365"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -2.15%  "

$ws.Range("D12").Value = "1.695.17"
$ws.Range("E12").Value = "  -5.70%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "4.521"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -5.43%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.5785"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -5.72%  "

$ws.Range("D15").Value = "1.925.76"
$ws.Range("E15").Value = "  -5.56%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.000008510"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -5.65%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "65.40"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -12.96%  "

$ws.Range("D18").Value = "26.624.28"
$ws.Range("E18").Value = "  -7.01%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "4.986"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -8.56%  "

$ws.Range("E20").Value = "  +0.21%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "10.97"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -4.37%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "186.63"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -10.83%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "6.256"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -8.12%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "1.007"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.18%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "144.51"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -5.57%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "7.467"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -6.97%  "

$ws.Range("E27").Value = "  -7.29%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "15.79"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -3.64%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.340"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -6.05%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.05737"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -5.86%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.339"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -5.79%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "3.522"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -6.74%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "3.509"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -7.75%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.641"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -5.02%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.020"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -2.40%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.5997"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -6.14%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.361"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -5.56%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "2.683"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -1.15%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.01619"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -4.52%  "

$ws.Range("D40").Value = "1.101.97"
$ws.Range("E40").Value = "  -3.50%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.8581"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -2.27%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "5.833"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -8.96%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.003"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.36%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "99.49"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.56%  "

$ws.Range("D45").Value = "1.852.28"
$ws.Range("E45").Value = "  -5.00%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.00000000118"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +6.51%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "56.53"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -5.59%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "1.005"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +0.68%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "8.093"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -3.12%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.4327"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -3.49%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.05236"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -4.49%  "
